$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E11").Value = 116000

$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

$ws.Range("C16").Value = "20300190"
$ws.Range("D16").Value = "KATHERIN MONTILLA PARALES"
$ws.Range("E16").Value = "2103"
$ws.Range("F16").Value = 116000
$ws.Range("G16").Value = 2900000

$ws.Range("17:22").Delete()
